# doc(./presentation): fixed mistypes. added pdf version of presentation
#
# 1) Fix typo "Лапытова" -> "Латыпова" in the "Проверил:" line on slide 1.
# 2) Refresh the cached "datetimeFigureOut" footer date (11/30/2021 -> 12/5/2021)
#    on the slide master and on every slide layout.

$p = $ppt.ActivePresentation

# --- 1) Fix the reviewer's name typo on slide 1 -----------------------------
$s1 = $p.Slides.Item(1)
for ($i = 1; $i -le $s1.Shapes.Count; $i++) {
    $sh = $s1.Shapes.Item($i)
    if ($sh.Name -eq "TextBox 4") {
        # Clearing first collapses the run back to a single plain run, so the
        # reassignment below produces one run per paragraph (matching how
        # PowerPoint normalizes a manually retyped line) instead of leaving
        # the old split ("Проверил: " / "Лапытова" / " О. А.") runs in place.
        $sh.TextFrame.TextRange.Text = ""
        $sh.TextFrame.TextRange.Text = "Выполнил: Воронцов С. А.`rПроверил: Латыпова О. А."
    }
}

# --- 2) Update the cached date field text everywhere it appears ------------
$oldDate = "11/30/2021"
$newDate = "12/5/2021"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*" -and $sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}
